$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column J ("Oficina (punto de entrega)") entirely, shifting all
# subsequent columns one position to the left.
$ws.Columns("J").Delete()
